# ============================================================
# Add a new '2022-Q1' sheet (between '2021-Q4' and the summary
# sheet '总计') with per-fund holding detail, and prepend a new
# '2022-Q1' row to the '总计' (total) summary sheet.
# ============================================================
$wb = $excel.ActiveWorkbook

# ---- 1) Insert new sheet "2022-Q1" after "2021-Q4" ----
$q4Sheet = $wb.Worksheets.Item(3)
$newWs = $wb.Worksheets.Add($null, $q4Sheet)
$newWs.Name = "2022-Q1"

# Force "Text" number format on the numeric-looking columns before
# writing their values, so fund codes / percentages / scale keep
# their exact original formatting (leading zeros, trailing zeros,
# fixed decimal places) instead of being auto-coerced to numbers.
$newWs.Range("B2:B48").NumberFormat = "@"
$newWs.Range("D2:D48").NumberFormat = "@"
$newWs.Range("E2:E48").NumberFormat = "@"
$newWs.Range("F2:F48").NumberFormat = "@"
$newWs.Range("G2:G47").NumberFormat = "@"

# ---- Build the full A1:H48 grid in one shot ----
$grid = New-Object "object[,]" 48,8

# Header row (row index 0 => sheet row 1). Column A header is blank.
$grid[0,1] = "基金代码"
$grid[0,2] = "基金名称"
$grid[0,3] = "基金规模"
$grid[0,4] = "股票总仓位"
$grid[0,5] = "仓位占比"
$grid[0,6] = "持有市值(亿元)"
$grid[0,7] = "仓位排名"

# Data rows (row index 1..47 => sheet rows 2..48)
$grid[1,0] = 0
$grid[1,1] = "510900"
$grid[1,2] = "易方达恒生国企(QDII-ETF)"
$grid[1,3] = "101.39"
$grid[1,4] = "97.12"
$grid[1,5] = "4.60"
$grid[1,6] = "4.6639"
$grid[1,7] = 5
$grid[2,0] = 1
$grid[2,1] = "159920"
$grid[2,2] = "华夏恒生ETF(QDII)"
$grid[2,3] = "151.31"
$grid[2,4] = "95.19"
$grid[2,5] = "2.77"
$grid[2,6] = "4.1913"
$grid[2,7] = 8
$grid[3,0] = 2
$grid[3,1] = "501025"
$grid[3,2] = "鹏华港股通中证香港银行投资指数（LOF）A"
$grid[3,3] = "9.81"
$grid[3,4] = "94.47"
$grid[3,5] = "14.05"
$grid[3,6] = "1.3783"
$grid[3,7] = 3
$grid[4,0] = 3
$grid[4,1] = "513550"
$grid[4,2] = "华泰柏瑞中证港股通50ETF"
$grid[4,3] = "31.28"
$grid[4,4] = "98.89"
$grid[4,5] = "3.87"
$grid[4,6] = "1.2105"
$grid[4,7] = 7
$grid[5,0] = 4
$grid[5,1] = "010365"
$grid[5,2] = "鹏华港股通中证香港银行投资指数（LOF）C"
$grid[5,3] = "6.07"
$grid[5,4] = "94.47"
$grid[5,5] = "14.05"
$grid[5,6] = "0.8528"
$grid[5,7] = 3
$grid[6,0] = 5
$grid[6,1] = "513660"
$grid[6,2] = "华夏沪港通恒生ETF"
$grid[6,3] = "19.61"
$grid[6,4] = "97.34"
$grid[6,5] = "3.11"
$grid[6,6] = "0.6099"
$grid[6,7] = 7
$grid[7,0] = 6
$grid[7,1] = "501050"
$grid[7,2] = "华夏沪港通上证50AH优选指数（LOF）A"
$grid[7,3] = "25.94"
$grid[7,4] = "92.28"
$grid[7,5] = "2.30"
$grid[7,6] = "0.5966"
$grid[7,7] = 10
$grid[8,0] = 7
$grid[8,1] = "006809"
$grid[8,2] = "泰康港股通中证香港银行投资指数A"
$grid[8,3] = "1.99"
$grid[8,4] = "94.73"
$grid[8,5] = "14.05"
$grid[8,6] = "0.2796"
$grid[8,7] = 3
$grid[9,0] = 8
$grid[9,1] = "159850"
$grid[9,2] = "华夏恒生中国企业ETF（QDII）"
$grid[9,3] = "6.11"
$grid[9,4] = "93.95"
$grid[9,5] = "4.43"
$grid[9,6] = "0.2707"
$grid[9,7] = 5
$grid[10,0] = 9
$grid[10,1] = "010010"
$grid[10,2] = "国投瑞银港股通6个月定期开放股票"
$grid[10,3] = "8.09"
$grid[10,4] = "93.58"
$grid[10,5] = "2.64"
$grid[10,6] = "0.2136"
$grid[10,7] = 9
$grid[11,0] = 10
$grid[11,1] = "159960"
$grid[11,2] = "平安港股通恒生中国企业ETF"
$grid[11,3] = "4.08"
$grid[11,4] = "96.86"
$grid[11,5] = "5.20"
$grid[11,6] = "0.2122"
$grid[11,7] = 4
$grid[12,0] = 11
$grid[12,1] = "501301"
$grid[12,2] = "华宝港股通恒生中国(香港上市)25指数(LOF)A"
$grid[12,3] = "2.82"
$grid[12,4] = "94.73"
$grid[12,5] = "6.78"
$grid[12,6] = "0.1912"
$grid[12,7] = 4
$grid[13,0] = 12
$grid[13,1] = "513600"
$grid[13,2] = "南方恒生ETF"
$grid[13,3] = "5.89"
$grid[13,4] = "99.00"
$grid[13,5] = "3.04"
$grid[13,6] = "0.1791"
$grid[13,7] = 7
$grid[14,0] = 13
$grid[14,1] = "159954"
$grid[14,2] = "南方恒生中国企业ETF"
$grid[14,3] = "3.35"
$grid[14,4] = "103.89"
$grid[14,5] = "4.89"
$grid[14,6] = "0.1638"
$grid[14,7] = 4
$grid[15,0] = 14
$grid[15,1] = "010204"
$grid[15,2] = "中银港股通优势成长股票"
$grid[15,3] = "3.19"
$grid[15,4] = "83.00"
$grid[15,5] = "4.89"
$grid[15,6] = "0.1560"
$grid[15,7] = 5
$grid[16,0] = 15
$grid[16,1] = "006810"
$grid[16,2] = "泰康港股通中证香港银行投资指数C"
$grid[16,3] = "0.90"
$grid[16,4] = "94.73"
$grid[16,5] = "14.05"
$grid[16,6] = "0.1264"
$grid[16,7] = 3
$grid[17,0] = 16
$grid[17,1] = "160717"
$grid[17,2] = "嘉实恒生中国企业指数(QDII-LOF)"
$grid[17,3] = "2.57"
$grid[17,4] = "94.76"
$grid[17,5] = "4.46"
$grid[17,6] = "0.1146"
$grid[17,7] = 5
$grid[18,0] = 17
$grid[18,1] = "161831"
$grid[18,2] = "银华恒生国企指数（QDII-LOF）"
$grid[18,3] = "2.29"
$grid[18,4] = "86.34"
$grid[18,5] = "4.70"
$grid[18,6] = "0.1076"
$grid[18,7] = 5
$grid[19,0] = 18
$grid[19,1] = "164705"
$grid[19,2] = "汇添富恒生指数（QDII-LOF）A"
$grid[19,3] = "2.96"
$grid[19,4] = "92.23"
$grid[19,5] = "2.69"
$grid[19,6] = "0.0796"
$grid[19,7] = 8
$grid[20,0] = 19
$grid[20,1] = "501310"
$grid[20,2] = "华宝标普沪港深中国增强价值指数（LOF）A"
$grid[20,3] = "1.44"
$grid[20,4] = "94.80"
$grid[20,5] = "4.91"
$grid[20,6] = "0.0707"
$grid[20,7] = 2
$grid[21,0] = 20
$grid[21,1] = "006355"
$grid[21,2] = "华宝港股通恒生中国(香港上市)25指数(LOF)C"
$grid[21,3] = "1.02"
$grid[21,4] = "94.73"
$grid[21,5] = "6.78"
$grid[21,6] = "0.0692"
$grid[21,7] = 4
$grid[22,0] = 21
$grid[22,1] = "517080"
$grid[22,2] = "汇添富中证沪港深500ETF"
$grid[22,3] = "5.68"
$grid[22,4] = "91.59"
$grid[22,5] = "0.96"
$grid[22,6] = "0.0545"
$grid[22,7] = 10
$grid[23,0] = 22
$grid[23,1] = "007107"
$grid[23,2] = "太平 MSCI 香港价值增强指数A"
$grid[23,3] = "1.05"
$grid[23,4] = "93.78"
$grid[23,5] = "5.07"
$grid[23,6] = "0.0532"
$grid[23,7] = 6
$grid[24,0] = 23
$grid[24,1] = "517100"
$grid[24,2] = "富国中证沪港深500ETF"
$grid[24,3] = "4.13"
$grid[24,4] = "99.22"
$grid[24,5] = "1.01"
$grid[24,6] = "0.0417"
$grid[24,7] = 10
$grid[25,0] = 24
$grid[25,1] = "159712"
$grid[25,2] = "国泰中证港股通50ETF"
$grid[25,3] = "0.85"
$grid[25,4] = "95.21"
$grid[25,5] = "4.63"
$grid[25,6] = "0.0394"
$grid[25,7] = 6
$grid[26,0] = 25
$grid[26,1] = "007751"
$grid[26,2] = "景顺长城中证沪港深红利成长低波动指数A"
$grid[26,3] = "0.83"
$grid[26,4] = "91.29"
$grid[26,5] = "2.66"
$grid[26,6] = "0.0221"
$grid[26,7] = 6
$grid[27,0] = 26
$grid[27,1] = "513990"
$grid[27,2] = "招商上证港股通ETF"
$grid[27,3] = "0.59"
$grid[27,4] = "96.48"
$grid[27,5] = "2.95"
$grid[27,6] = "0.0174"
$grid[27,7] = 7
$grid[28,0] = 27
$grid[28,1] = "513680"
$grid[28,2] = "建信港股通恒生中国企业ETF"
$grid[28,3] = "0.28"
$grid[28,4] = "96.77"
$grid[28,5] = "5.57"
$grid[28,6] = "0.0156"
$grid[28,7] = 4
$grid[29,0] = 28
$grid[29,1] = "008407"
$grid[29,2] = "恒生前海恒生沪深港通细分行业龙头指数A"
$grid[29,3] = "0.37"
$grid[29,4] = "93.40"
$grid[29,5] = "4.01"
$grid[29,6] = "0.0148"
$grid[29,7] = 5
$grid[30,0] = 29
$grid[30,1] = "010789"
$grid[30,2] = "汇添富恒生指数（QDII-LOF）C"
$grid[30,3] = "0.37"
$grid[30,4] = "92.23"
$grid[30,5] = "2.69"
$grid[30,6] = "0.0100"
$grid[30,7] = 8
$grid[31,0] = 30
$grid[31,1] = "006658"
$grid[31,2] = "财通中证香港红利等权投资指数A"
$grid[31,3] = "0.20"
$grid[31,4] = "90.59"
$grid[31,5] = "4.89"
$grid[31,6] = "0.0098"
$grid[31,7] = 2
$grid[32,0] = 31
$grid[32,1] = "001942"
$grid[32,2] = "前海开源沪港深汇鑫灵活配置混合A"
$grid[32,3] = "0.10"
$grid[32,4] = "90.39"
$grid[32,5] = "8.20"
$grid[32,6] = "0.0082"
$grid[32,7] = 2
$grid[33,0] = 32
$grid[33,1] = "501309"
$grid[33,2] = "国泰恒生港股通指数（LOF）"
$grid[33,3] = "0.36"
$grid[33,4] = "92.35"
$grid[33,5] = "1.87"
$grid[33,6] = "0.0067"
$grid[33,7] = 8
$grid[34,0] = 33
$grid[34,1] = "001943"
$grid[34,2] = "前海开源沪港深汇鑫灵活配置混合C"
$grid[34,3] = "0.08"
$grid[34,4] = "90.39"
$grid[34,5] = "8.20"
$grid[34,6] = "0.0066"
$grid[34,7] = 2
$grid[35,0] = 34
$grid[35,1] = "160925"
$grid[35,2] = "大成中华沪深港300指数（LOF）A"
$grid[35,3] = "0.54"
$grid[35,4] = "93.14"
$grid[35,5] = "1.11"
$grid[35,6] = "0.0060"
$grid[35,7] = 10
$grid[36,0] = 35
$grid[36,1] = "166402"
$grid[36,2] = "浦银安盛中证锐联沪港深基本面100指数（LOF）"
$grid[36,3] = "0.19"
$grid[36,4] = "90.95"
$grid[36,5] = "3.03"
$grid[36,6] = "0.0058"
$grid[36,7] = 5
$grid[37,0] = 36
$grid[37,1] = "006395"
$grid[37,2] = "华夏沪港通上证50AH优选指数（LOF）C"
$grid[37,3] = "0.25"
$grid[37,4] = "92.28"
$grid[37,5] = "2.30"
$grid[37,6] = "0.0058"
$grid[37,7] = 10
$grid[38,0] = 37
$grid[38,1] = "517010"
$grid[38,2] = "易方达中证沪港深500交易型开放式指数证券投资基金"
$grid[38,3] = "0.44"
$grid[38,4] = "91.01"
$grid[38,5] = "1.06"
$grid[38,6] = "0.0047"
$grid[38,7] = 10
$grid[39,0] = 38
$grid[39,1] = "007397"
$grid[39,2] = "华宝标普沪港深中国增强价值指数（LOF）C"
$grid[39,3] = "0.09"
$grid[39,4] = "94.80"
$grid[39,5] = "4.91"
$grid[39,6] = "0.0044"
$grid[39,7] = 2
$grid[40,0] = 39
$grid[40,1] = "011647"
$grid[40,2] = "博时港股通红利精选混合A"
$grid[40,3] = "0.13"
$grid[40,4] = "92.10"
$grid[40,5] = "3.23"
$grid[40,6] = "0.0042"
$grid[40,7] = 9
$grid[41,0] = 40
$grid[41,1] = "517170"
$grid[41,2] = "华夏中证沪港深500交易型开放式指数证券投资基金"
$grid[41,3] = "0.33"
$grid[41,4] = "94.49"
$grid[41,5] = "0.97"
$grid[41,6] = "0.0032"
$grid[41,7] = 10
$grid[42,0] = 41
$grid[42,1] = "008408"
$grid[42,2] = "恒生前海恒生沪深港通细分行业龙头指数C"
$grid[42,3] = "0.08"
$grid[42,4] = "93.40"
$grid[42,5] = "4.01"
$grid[42,6] = "0.0032"
$grid[42,7] = 5
$grid[43,0] = 42
$grid[43,1] = "006659"
$grid[43,2] = "财通中证香港红利等权投资指数C"
$grid[43,3] = "0.05"
$grid[43,4] = "90.59"
$grid[43,5] = "4.89"
$grid[43,6] = "0.0024"
$grid[43,7] = 2
$grid[44,0] = 43
$grid[44,1] = "007760"
$grid[44,2] = "景顺长城中证沪港深红利成长低波动指数C"
$grid[44,3] = "0.06"
$grid[44,4] = "91.29"
$grid[44,5] = "2.66"
$grid[44,6] = "0.0016"
$grid[44,7] = 6
$grid[45,0] = 44
$grid[45,1] = "011648"
$grid[45,2] = "博时港股通红利精选混合C"
$grid[45,3] = "0.02"
$grid[45,4] = "92.10"
$grid[45,5] = "3.23"
$grid[45,6] = "0.0006"
$grid[45,7] = 9
$grid[46,0] = 45
$grid[46,1] = "008973"
$grid[46,2] = "大成中华沪深港300指数(LOF)C"
$grid[46,3] = "0.02"
$grid[46,4] = "93.14"
$grid[46,5] = "1.11"
$grid[46,6] = "0.0002"
$grid[46,7] = 10
$grid[47,0] = 46
$grid[47,1] = "007108"
$grid[47,2] = "太平 MSCI 香港价值增强指数C"
$grid[47,3] = "0.00"
$grid[47,4] = "93.78"
$grid[47,5] = "5.07"
$grid[47,6] = 0
$grid[47,7] = 6

$newWs.Range("A1:H48").Value = $grid

# ---- Styling: copy the header-row / index-column look from the
#      "2021-Q4" sheet (bold, centered, bordered) onto the new sheet ----
$q4Sheet.Range("B1:H1").Copy()
$newWs.Range("B1:H1").PasteSpecial(-4122)
$q4Sheet.Range("A2").Copy()
$newWs.Range("A2:A48").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- 2) Prepend a "2022-Q1" row to the "总计" (total) sheet ----
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Rows.Item(2).Insert()
$totalWs.Range("B2:D2").ClearFormats()
$totalWs.Range("A3").Copy()
$totalWs.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$totalWs.Range("A2").Value = 0
$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 47
$totalWs.Range("D2").Value = 16.08

# Renumber the existing rows index column (0,1,2,3,...)
$totalWs.Range("A3").Value = 1
$totalWs.Range("A4").Value = 2
$totalWs.Range("A5").Value = 3

